# Developer guide diagram edits — HighLevelSequenceDiagrams
# (see commit: "Editing the developer guide")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "delete 1" -> "delete d1" -------------------------------------------

# TextBox 23 : the plain label on the actor-lifeline arrow
$s.Shapes.Item("TextBox 23").TextFrame.TextRange.Text = "delete d1"

# TextBox 25 : execute(“delete 1”) -> execute(“delete d1”), plus the shape
# grows a little (and shifts left) to keep the new, longer label centered
# on its arrow. Resize *before* changing the text so the auto-fit height
# isn't thrown off by a momentary width/text mismatch.
$execShape = $s.Shapes.Item("TextBox 25")
$execShape.Left = (2084433 + 0.5) / 12700
$execShape.Width = (1506585 + 0.5) / 12700
$execShape.TextFrame.TextRange.Text = "execute(" + [char]8220 + "delete d1" + [char]8221 + ")"

# --- method / event renames ------------------------------------------------

# TextBox 28 : deletePerson(p) -> deleteDeadline(p)  (only the method-name
# run changes; the "(p)" run is untouched)
$s.Shapes.Item("TextBox 28").TextFrame.TextRange.Characters(1, 12).Text = "deleteDeadline"

# TextBox 32 & TextBox 61 : post(AddressBookChangedEvent) -> post(ToDoListChangedEvent)
$s.Shapes.Item("TextBox 32").TextFrame.TextRange.Characters(6, 23).Text = "ToDoListChangedEvent"
$s.Shapes.Item("TextBox 61").TextFrame.TextRange.Characters(6, 23).Text = "ToDoListChangedEvent"

# TextBox 73 & TextBox 49 : handleAddresssBookChangedEvent() -> handleToDoListChangedEvent()
$s.Shapes.Item("TextBox 73").TextFrame.TextRange.Characters(1, 30).Text = "handleToDoListChangedEvent"
$s.Shapes.Item("TextBox 49").TextFrame.TextRange.Characters(1, 30).Text = "handleToDoListChangedEvent"

# --- presentation-level slide guides ---------------------------------------
# Two new guides shown on every slide: a horizontal one at 186pt and a
# vertical one at 360pt (PowerPoint's standard custom-guide grey).
$null = $p.Guides.Add(1, 186)
$null = $p.Guides.Add(2, 360)
